$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header: volume number and week-covering date range
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# Update weekly crime statistics table (rows 14-30)
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 300
$ws.Range("I14").Value = 44
$ws.Range("K14").Value = -2.222222222222
$ws.Range("L14").Value = -25.423728813559
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -84.722222222222
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -42.857142857142
$ws.Range("G15").Value = 22
$ws.Range("H15").Value = -18.181818181818
$ws.Range("I15").Value = 142
$ws.Range("J15").Value = 148
$ws.Range("K15").Value = -4.054054054054
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 7.575757575757
$ws.Range("N15").Value = -62.532981530343
$ws.Range("C16").Value = 49
$ws.Range("D16").Value = 48
$ws.Range("E16").Value = 2.083333333333
$ws.Range("F16").Value = 209
$ws.Range("G16").Value = 206
$ws.Range("H16").Value = 1.456310679611
$ws.Range("I16").Value = 1484
$ws.Range("J16").Value = 1592
$ws.Range("K16").Value = -6.783919597989
$ws.Range("L16").Value = 23.563696919234
$ws.Range("M16").Value = -28.756601056169
$ws.Range("N16").Value = -85.155546663999
$ws.Range("C17").Value = 63
$ws.Range("D17").Value = 87
$ws.Range("E17").Value = -27.586206896551
$ws.Range("F17").Value = 333
$ws.Range("G17").Value = 377
$ws.Range("H17").Value = -11.671087533156
$ws.Range("I17").Value = 2625
$ws.Range("J17").Value = 2577
$ws.Range("K17").Value = 1.862630966239
$ws.Range("L17").Value = 26.080691642651
$ws.Range("M17").Value = 26.445086705202
$ws.Range("N17").Value = -50.274673233567
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 47
$ws.Range("E18").Value = -23.404255319148
$ws.Range("F18").Value = 181
$ws.Range("G18").Value = 193
$ws.Range("H18").Value = -6.217616580310
$ws.Range("I18").Value = 1271
$ws.Range("J18").Value = 1483
$ws.Range("K18").Value = -14.295347269049
$ws.Range("L18").Value = 8.818493150684
$ws.Range("M18").Value = -32.068412613575
$ws.Range("N18").Value = -82.705129949653
$ws.Range("C19").Value = 96
$ws.Range("D19").Value = 124
$ws.Range("E19").Value = -22.580645161290
$ws.Range("F19").Value = 469
$ws.Range("G19").Value = 515
$ws.Range("H19").Value = -8.932038834951
$ws.Range("I19").Value = 3558
$ws.Range("J19").Value = 3611
$ws.Range("K19").Value = -1.467737468845
$ws.Range("L19").Value = 34.061793519216
$ws.Range("M19").Value = 42.206235011990
$ws.Range("N19").Value = -13.282963685108
$ws.Range("C20").Value = 44
$ws.Range("D20").Value = 31
$ws.Range("E20").Value = 41.935483870967
$ws.Range("F20").Value = 173
$ws.Range("G20").Value = 151
$ws.Range("H20").Value = 14.569536423841
$ws.Range("I20").Value = 1111
$ws.Range("J20").Value = 1087
$ws.Range("K20").Value = 2.207911683532
$ws.Range("L20").Value = 23.170731707317
$ws.Range("M20").Value = 27.554535017221
$ws.Range("N20").Value = -80.805114029025
$ws.Range("C21").Value = 295
$ws.Range("D21").Value = 344
$ws.Range("E21").Value = -14.244186046511
$ws.Range("F21").Value = 1391
$ws.Range("G21").Value = 1466
$ws.Range("H21").Value = -5.115961800818
$ws.Range("I21").Value = 10235
$ws.Range("J21").Value = 10543
$ws.Range("K21").Value = -2.921369629137
$ws.Range("L21").Value = 24.695419103313
$ws.Range("M21").Value = 6.359763067650
$ws.Range("N21").Value = -69.155893077780
$ws.Range("C22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 12.5
$ws.Range("I22").Value = 177
$ws.Range("J22").Value = 219
$ws.Range("K22").Value = -19.178082191780
$ws.Range("L22").Value = 23.776223776223
$ws.Range("M22").Value = -31.128404669260
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 140
$ws.Range("H23").Value = -20.714285714285
$ws.Range("I23").Value = 986
$ws.Range("J23").Value = 952
$ws.Range("K23").Value = 3.571428571428
$ws.Range("L23").Value = 15.052508751458
$ws.Range("M23").Value = 41.463414634146
$ws.Range("C24").Value = 261
$ws.Range("D24").Value = 249
$ws.Range("E24").Value = 4.819277108433
$ws.Range("F24").Value = 1053
$ws.Range("G24").Value = 1119
$ws.Range("H24").Value = -5.898123324396
$ws.Range("I24").Value = 7707
$ws.Range("J24").Value = 8086
$ws.Range("K24").Value = -4.687113529557
$ws.Range("L24").Value = 25.869671729544
$ws.Range("M24").Value = 22.722929936305
$ws.Range("C25").Value = 125
$ws.Range("D25").Value = 105
$ws.Range("E25").Value = 19.047619047619
$ws.Range("F25").Value = 493
$ws.Range("G25").Value = 437
$ws.Range("H25").Value = 12.814645308924
$ws.Range("I25").Value = 3807
$ws.Range("J25").Value = 3739
$ws.Range("K25").Value = 1.818668093073
$ws.Range("L25").Value = 36.110117983553
$ws.Range("M25").Value = -23.106443142799
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -12.121212121212
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 230
$ws.Range("K26").Value = -4.347826086956
$ws.Range("L26").Value = -7.172995780590
$ws.Range("C27").Value = 16
$ws.Range("D27").Value = 17
$ws.Range("E27").Value = -5.882352941176
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 52
$ws.Range("H27").Value = 3.846153846153
$ws.Range("I27").Value = 383
$ws.Range("J27").Value = 391
$ws.Range("K27").Value = -2.046035805626
$ws.Range("L27").Value = -10.304449648711
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = -51.219512195122
$ws.Range("I28").Value = 146
$ws.Range("J28").Value = 221
$ws.Range("K28").Value = -33.936651583710
$ws.Range("L28").Value = -42.063492063492
$ws.Range("M28").Value = -54.938271604938
$ws.Range("N28").Value = -87.700084245998
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 17
$ws.Range("G29").Value = 33
$ws.Range("H29").Value = -48.484848484848
$ws.Range("I29").Value = 127
$ws.Range("J29").Value = 185
$ws.Range("K29").Value = -31.351351351351
$ws.Range("L29").Value = -39.234449760765
$ws.Range("M29").Value = -51.340996168582
$ws.Range("N29").Value = -88.086303939962
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 46
$ws.Range("K30").Value = -13.043478260869
